$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.190.20"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "1.850.55"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7048"
$ws.Range("E5").Value = "  -4.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.99"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3054"
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07439"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.44"
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08151"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").Value = "1.889.96"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7289"
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.18"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("D16").Value = "29.415.57"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.783"
$ws.Range("E17").Value = "  -6.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.97"
$ws.Range("E18").Value = "  -4.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007665"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "2.140.57"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.604"
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.025"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1462"
$ws.Range("E26").Value = "  -6.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.83"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.14"
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.977"
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.410"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.515"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.496"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.012"
$ws.Range("E33").Value = "  -4.58%  "
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.189"
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.039"
$ws.Range("E36").Value = "  +3.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7078"
$ws.Range("E37").Value = "  -8.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.659"
$ws.Range("E39").Value = "  -4.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.679"
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9431"
$ws.Range("E41").Value = "  +7.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.021"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4309"
$ws.Range("E43").Value = "  -5.84%  "
$ws.Range("D44").Value = "1.067.54"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.50"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.62"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "2.039.02"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.749"
$ws.Range("E49").Value = "  -5.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.053"
$ws.Range("E50").Value = "  -7.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.109"
$ws.Range("E51").Value = "  -4.67%  "
